$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.061121362606882
$ws.Range("D2").Value = 1.060064216073359
$ws.Range("E2").Value = 1.066089893537598
$ws.Range("F2").Value = 1.075132539211116
$ws.Range("I2").Value = 1.047037053148322
$ws.Range("J2").Value = 1.066098679675674
$ws.Range("K2").Value = 1.062792013712198
$ws.Range("L2").Value = 1.068801378960385
$ws.Range("M2").Value = 1.077819915998275
$ws.Range("N2").Value = 1.067612662084041
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.062622173286897
$ws.Range("D3").Value = 1.061402461530859
$ws.Range("E3").Value = 1.067429868840806
$ws.Range("F3").Value = 1.076564025129314
$ws.Range("I3").Value = 1.047423348009112
$ws.Range("J3").Value = 1.06725077828038
$ws.Range("K3").Value = 1.063943287993346
$ws.Range("L3").Value = 1.069955570550123
$ws.Range("M3").Value = 1.079067153536791
$ws.Range("N3").Value = 1.068766396800914
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.063591789410865
$ws.Range("D4").Value = 1.062267181350943
$ws.Range("E4").Value = 1.068295775496896
$ws.Range("F4").Value = 1.077489203470759
$ws.Range("I4").Value = 1.047670885645462
$ws.Range("J4").Value = 1.067994331806291
$ws.Range("K4").Value = 1.06468648491095
$ws.Range("L4").Value = 1.070700720504051
$ws.Range("M4").Value = 1.079872583911924
$ws.Range("N4").Value = 1.069511006258118
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.063999062165444
$ws.Range("D5").Value = 1.06263042462949
$ws.Range("E5").Value = 1.068659533821859
$ws.Range("F5").Value = 1.077877893899544
$ws.Range("I5").Value = 1.047774373034264
$ws.Range("J5").Value = 1.068306464781961
$ws.Range("K5").Value = 1.064998510204386
$ws.Range("L5").Value = 1.071013582456499
$ws.Range("M5").Value = 1.080210805619908
$ws.Range("N5").Value = 1.069823582498415
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.064067424514541
$ws.Range("D6").Value = 1.062691398233939
$ws.Range("E6").Value = 1.068720594847244
$ws.Range("F6").Value = 1.077943141982932
$ws.Range("I6").Value = 1.047791715223369
$ws.Range("J6").Value = 1.068358846644831
$ws.Range("K6").Value = 1.065050876461188
$ws.Range("L6").Value = 1.071066090105807
$ws.Range("M6").Value = 1.080267572339468
$ws.Range("N6").Value = 1.069876038749539
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.063597232790727
$ws.Range("D7").Value = 1.062272036137823
$ws.Range("E7").Value = 1.068300637102661
$ws.Range("F7").Value = 1.077494398163938
$ws.Range("I7").Value = 1.0476722707153
$ws.Range("J7").Value = 1.067998504331102
$ws.Range("K7").Value = 1.06469065583105
$ws.Range("L7").Value = 1.07070490254237
$ws.Range("M7").Value = 1.079877104737993
$ws.Range("N7").Value = 1.069515184708393
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.061628885423549
$ws.Range("D8").Value = 1.060516737367059
$ws.Range("E8").Value = 1.066542985425058
$ws.Range("F8").Value = 1.075616545703192
$ws.Range("I8").Value = 1.047168105883795
$ws.Range("J8").Value = 1.066488439964282
$ws.Range("K8").Value = 1.063181458759817
$ws.Range("L8").Value = 1.069191796371811
$ws.Range("M8").Value = 1.078241763679013
$ws.Range("N8").Value = 1.068002975876966
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.058148521043509
$ws.Range("D9").Value = 1.05741411358117
$ws.Range("E9").Value = 1.063436735351581
$ws.Range("F9").Value = 1.072298916246727
$ws.Range("I9").Value = 1.046261064035934
$ws.Range("J9").Value = 1.063812477495669
$ws.Range("K9").Value = 1.060508382497235
$ws.Range("L9").Value = 1.06651233407578
$ws.Range("M9").Value = 1.075347458038679
$ws.Range("N9").Value = 1.06532321323471
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.055819824135505
$ws.Range("D10").Value = 1.055338899685063
$ws.Range("E10").Value = 1.06135945736986
$ws.Range("F10").Value = 1.070080985109392
$ws.Range("I10").Value = 1.045643699757524
$ws.Range("J10").Value = 1.062018048560072
$ws.Range("K10").Value = 1.058716796369404
$ws.Range("L10").Value = 1.064716830942759
$ws.Range("M10").Value = 1.073409078103526
$ws.Range("N10").Value = 1.0635262360042
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.054809363712354
$ws.Range("D11").Value = 1.054438614037961
$ws.Range("E11").Value = 1.060458362666399
$ws.Range("F11").Value = 1.06911904417341
$ws.Range("I11").Value = 1.045373337852236
$ws.Range("J11").Value = 1.061238484066917
$ws.Range("K11").Value = 1.057938682044711
$ws.Range("L11").Value = 1.063937101836163
$ws.Range("M11").Value = 1.072567561245846
$ws.Range("N11").Value = 1.062745564440048
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.054433706291053
$ws.Range("D12").Value = 1.054103944589057
$ws.Range("E12").Value = 1.060123405501494
$ws.Range("F12").Value = 1.068761494209459
$ws.Range("I12").Value = 1.045272454036014
$ws.Range("J12").Value = 1.060948526986603
$ws.Range("K12").Value = 1.057649296673395
$ws.Range("L12").Value = 1.063647128931388
$ws.Range("M12").Value = 1.072254649112118
$ws.Range("N12").Value = 1.062455195587421
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.054514301045718
$ws.Range("D13").Value = 1.054175744352134
$ws.Range("E13").Value = 1.060195266399988
$ws.Range("F13").Value = 1.068838200976038
$ws.Range("I13").Value = 1.045294114789238
$ws.Range("J13").Value = 1.061010741618001
$ws.Range("K13").Value = 1.057711387173698
$ws.Range("L13").Value = 1.063709344903581
$ws.Range("M13").Value = 1.072321785112423
$ws.Range("N13").Value = 1.062517498570732
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.054778318484667
$ws.Range("D14").Value = 1.05441095556959
$ws.Range("E14").Value = 1.060430680153693
$ws.Range("F14").Value = 1.06908949396815
$ws.Range("I14").Value = 1.04536500815162
$ws.Range("J14").Value = 1.06121452416263
$ws.Range("K14").Value = 1.057914768723133
$ws.Range("L14").Value = 1.063913139698303
$ws.Range("M14").Value = 1.07254170271177
$ws.Range("N14").Value = 1.062721570509948
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.054940944683583
$ws.Range("D15").Value = 1.054555841860299
$ws.Range("E15").Value = 1.060575692932598
$ws.Range("F15").Value = 1.069244291543392
$ws.Range("I15").Value = 1.045408626937869
$ws.Range("J15").Value = 1.061340029151282
$ws.Range("K15").Value = 1.058040031028653
$ws.Range("L15").Value = 1.064038658245173
$ws.Range("M15").Value = 1.072677156562978
$ws.Range("N15").Value = 1.062847253730079
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.055886839645724
$ws.Range("D16").Value = 1.055398612155364
$ws.Range("E16").Value = 1.061419225325812
$ws.Range("F16").Value = 1.070144792403688
$ws.Range("I16").Value = 1.045661578520093
$ws.Range("J16").Value = 1.062069731071781
$ws.Range("K16").Value = 1.058768387266771
$ws.Range("L16").Value = 1.064768530725431
$ws.Range("M16").Value = 1.073464880166926
$ws.Range("N16").Value = 1.063577991911004
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.056479600817102
$ws.Range("D17").Value = 1.055926798319537
$ws.Range("E17").Value = 1.061947912266601
$ws.Range("F17").Value = 1.070709229247322
$ws.Range("I17").Value = 1.045819432760139
$ws.Range("J17").Value = 1.062526762345547
$ws.Range("K17").Value = 1.059224633201579
$ws.Range("L17").Value = 1.065225749564203
$ws.Range("M17").Value = 1.073958408348848
$ws.Range("N17").Value = 1.0640356722216
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.056825144512546
$ws.Range("D18").Value = 1.056234716420431
$ws.Range("E18").Value = 1.062256130629674
$ws.Range("F18").Value = 1.071038305502918
$ws.Range("I18").Value = 1.045911213428143
$ws.Range("J18").Value = 1.062793093937972
$ws.Range("K18").Value = 1.059490527861984
$ws.Range("L18").Value = 1.06549221954437
$ws.Range("M18").Value = 1.074246064586287
$ws.Range("N18").Value = 1.064302382035435
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.056942931675463
$ws.Range("D19").Value = 1.056339680943486
$ws.Range("E19").Value = 1.062361198917423
$ws.Range("F19").Value = 1.071150486752521
$ws.Range("I19").Value = 1.045942458642919
$ws.Range("J19").Value = 1.062883864433542
$ws.Range("K19").Value = 1.059581152957715
$ws.Range("L19").Value = 1.065583042137712
$ws.Range("M19").Value = 1.074344112450922
$ws.Range("N19").Value = 1.064393281435521
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.056416024305614
$ws.Range("D20").Value = 1.05587014595353
$ws.Range("E20").Value = 1.061891205326626
$ws.Range("F20").Value = 1.070648686107225
$ws.Range("I20").Value = 1.045802526824685
$ws.Range("J20").Value = 1.062477752814376
$ws.Range("K20").Value = 1.05917570574209
$ws.Range("L20").Value = 1.065176716909879
$ws.Range("M20").Value = 1.073905479258807
$ws.Range("N20").Value = 1.063986593091274
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.054700581040878
$ws.Range("D21").Value = 1.054341699054449
$ws.Range("E21").Value = 1.060361363642607
$ws.Range("F21").Value = 1.0690155011943
$ws.Range("I21").Value = 1.045344144526576
$ws.Range("J21").Value = 1.061154526161369
$ws.Range("K21").Value = 1.057854887891615
$ws.Range("L21").Value = 1.063853136837249
$ws.Range("M21").Value = 1.072476951772199
$ws.Range("N21").Value = 1.062661487304643
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.053620114845466
$ws.Range("D22").Value = 1.053379176983442
$ws.Range("E22").Value = 1.059398038731478
$ws.Range("F22").Value = 1.067987247520459
$ws.Range("I22").Value = 1.045053282232815
$ws.Range("J22").Value = 1.060320288982568
$ws.Range("K22").Value = 1.057022356510156
$ws.Range("L22").Value = 1.063018939745399
$ws.Range("M22").Value = 1.071576836392377
$ws.Range("N22").Value = 1.06182606541336
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.054193073447408
$ws.Range("D23").Value = 1.05388957529584
$ws.Range("E23").Value = 1.059908855593108
$ws.Range("F23").Value = 1.068532480018964
$ws.Range("I23").Value = 1.04520772682861
$ws.Range("J23").Value = 1.060762751533884
$ws.Range("K23").Value = 1.057463896614422
$ws.Range("L23").Value = 1.063461356125571
$ws.Range("M23").Value = 1.072054191137048
$ws.Range("N23").Value = 1.06226915631225
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.056444752424712
$ws.Range("D24").Value = 1.05589574522603
$ws.Range("E24").Value = 1.06191682923352
$ws.Range("F24").Value = 1.070676043410605
$ws.Range("I24").Value = 1.045810166794799
$ws.Range("J24").Value = 1.062499898875419
$ws.Range("K24").Value = 1.059197814653143
$ws.Range("L24").Value = 1.065198873329725
$ws.Range("M24").Value = 1.073929396286675
$ws.Range("N24").Value = 1.064008770602264
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.0590497345611
$ws.Range("D25").Value = 1.058217387221055
$ws.Range("E25").Value = 1.064240884469442
$ws.Range("F25").Value = 1.07315766337813
$ws.Range("I25").Value = 1.046497779019559
$ws.Range("J25").Value = 1.064506095090791
$ws.Range("K25").Value = 1.061201092440877
$ws.Range("L25").Value = 1.067206635169303
$ws.Range("M25").Value = 1.076097237910096
$ws.Range("N25").Value = 1.066017815846376
